# BOT; UPDATE DATA
# Inserts three new daily rows (2020-04-08 .. 2020-04-10) into the "相談件数"
# sheet just above the trailing "filler" / footnote rows, fills in their
# values, and updates the sheet's print area / dimension / view state to
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Insert 3 new rows above the old row 74 (the "filler" row), pushing the
#     existing filler (74) and footnote (75) rows down to 77/78. Inserting
#     (rather than just writing past the end) lets the new rows inherit the
#     existing number formats / styles from the surrounding data rows. ---
$ws.Range("A74:E76").Insert(-4121) | Out-Null

# --- Fill in the newly inserted rows with the new daily figures. ---
$ws.Cells.Item(74, 1).Value = 43929
$ws.Cells.Item(74, 2).Value = 938
$ws.Cells.Item(74, 3).Value = 17734
$ws.Cells.Item(74, 4).Value = 158
$ws.Cells.Item(74, 5).Value = 4546

$ws.Cells.Item(75, 1).Value = 43930
$ws.Cells.Item(75, 2).Value = 892
$ws.Cells.Item(75, 3).Value = 18626
$ws.Cells.Item(75, 4).Value = 171
$ws.Cells.Item(75, 5).Value = 4735

$ws.Cells.Item(76, 1).Value = 43931
$ws.Cells.Item(76, 2).Value = 926
$ws.Cells.Item(76, 3).Value = 19552
$ws.Cells.Item(76, 4).Value = 137
$ws.Cells.Item(76, 5).Value = 4872

# --- Update the named print area so it covers the three extra rows
#     (A1:E79 -> A1:E82). ---
foreach ($n in $wb.Names) {
    $n.RefersTo = "=相談件数!`$A`$1:`$E`$82"
}

# --- Re-establish the frozen header row/column and scroll the view down so
#     the newly-entered rows are visible, then leave the selection where the
#     user would naturally continue entering data (first empty cell under
#     the footnote row). ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 69
$win.ScrollColumn = 3
$ws.Range("B80").Select() | Out-Null
